# Pooh Points: normal 20260129
# Rename owner/team short-codes to their full team names on both sheets,
# and widen the owner columns to fit the new, longer names.

$wb = $excel.ActiveWorkbook

$playersSheet = $wb.Worksheets.Item("Players")
$totalsSheet  = $wb.Worksheets.Item("OwnerTotals")

# Old short-code -> new full name
$nameMap = @{
    "Booz" = "Boozers Losers"
    "CDL"  = "The Backslashers"
    "Clay" = "Hilton Heads"
    "Hal"  = "Three Dawg Nite"
    "Mark" = "Bend Rimmers"
    "Ron"  = "G-Flop"
    "Tar"  = "The Oddities"
}

# --- "Players" sheet: column B ("owner") holds the short codes, rows 2-38 ---
$lastRow = $playersSheet.Cells.Item($playersSheet.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 38 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $playersSheet.Cells.Item($r, 2)
    $current = $cell.Value2
    if ($nameMap.ContainsKey($current)) {
        $cell.Value = $nameMap[$current]
    }
}

# --- "OwnerTotals" sheet: column A ("owner") holds the short codes, rows 2-8 ---
$lastRowTotals = $totalsSheet.Cells.Item($totalsSheet.Rows.Count, 1).End(-4162).Row
if ($lastRowTotals -lt 2) { $lastRowTotals = 8 }

for ($r = 2; $r -le $lastRowTotals; $r++) {
    $cell = $totalsSheet.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($nameMap.ContainsKey($current)) {
        $cell.Value = $nameMap[$current]
    }
}

# --- Widen columns to fit the longer owner/team names ---
# Stored column <col width="..."> is 5/6 wider than the ColumnWidth value we set,
# so subtract 5/6 to land exactly on the target stored width of 18.
$targetWidth = 18 - (5 / 6)

$playersSheet.Columns.Item(2).ColumnWidth = $targetWidth   # "Players" sheet, column B (owner)
$totalsSheet.Columns.Item(1).ColumnWidth  = $targetWidth   # "OwnerTotals" sheet, column A (owner)
